# Insert a new weekly price record for "Haba" (Vega Central Mapocho de
# Santiago) right before the existing row 313. Excel shifts every row
# from 313 downward to 314.. (old row 414 becomes row 415), which is
# exactly the behaviour captured by the diff (dimension A1:R414 -> A1:R415).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 313:414 down one row, leaving a blank row 313 behind.
$ws.Rows.Item(313).Insert()

# Fill in the new row 313 with the new record's data.
$ws.Range("A313").Value = 9
$ws.Range("B313").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C313").Value = "Metropolitana"
$ws.Range("D313").Value = 45229
$ws.Range("E313").Value = 13
$ws.Range("F313").Value = 100112026
$ws.Range("G313").Value = "Haba"
$ws.Range("H313").Value = "Sin especificar"
$ws.Range("I313").Value = "Primera"
$ws.Range("J313").Value = 70
$ws.Range("K313").Value = 6000
$ws.Range("L313").Value = 7000
$ws.Range("M313").Value = 6500
$ws.Range("N313").Value = "`$/caja 20 kilos"
$ws.Range("O313").Value = "Provincia de Melipilla"
$ws.Range("P313").Value = 325
$ws.Range("Q313").Value = 20
$ws.Range("R313").Value = "Hortaliza"
